# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing cell counts (E, K) and the derived
# expression/specificity metrics (G, H, M, N, O, P, Q, R, S, T) for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3
        G = 2.520808
        H = 7.562424
        K = 3
        M = 2.553279333333334
        N = 7.659838000000001
        O = 0.1645043904057808
        P = 0.1645043904057808
        Q = 6.436326969701335
        R = 57.926942727312
        S = 0.1645043904057808
        T = 0.1645043904057808
    }
    3 = @{
        E = 3
        G = 2.520808
        H = 7.562424
        K = 3
        M = 8.058662
        N = 24.175986
        O = 0.5192088709172035
        P = 0.5192088709172035
        Q = 20.314339638896
        R = 182.829056750064
        S = 0.5192088709172035
        T = 0.5192088709172035
    }
    4 = @{
        E = 3
        G = 2.520808
        H = 7.562424
        K = 3
        M = 4.909099333333334
        N = 14.727298
        O = 0.3162867386770157
        P = 0.3162867386770157
        Q = 12.37489687226133
        R = 111.374071850352
        S = 0.3162867386770157
        T = 0.3162867386770157
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
